$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the remaining cells of row 54 (C54:G54) to log a new entry.
$ws.Range("C54").Value = "5:10PM"
$ws.Range("D54").Value = 20
$ws.Range("E54").Value = 170
$ws.Range("F54").Value = "Finish newProd page ui and code; start on newOrder page and code"
$ws.Range("G54").Value = "currently finished newProd page and code; majority of fields in newOrder page display correctly now. Next is to code the adding of products to the current order and successfully add the order to the db"

# Update the active selection on the sheet to G54, matching the saved view state.
$ws.Range("G54").Select()
